$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testCitizen")

# Column A (rows 1-8): "ulkeleris.." -> "ulkemis.."
$ws.Range("A1").Value = "ulkemis11"
$ws.Range("A2").Value = "ulkemis22"
$ws.Range("A3").Value = "ulkemis33"
$ws.Range("A4").Value = "ulkemis44"
$ws.Range("A5").Value = "ulkemis55"
$ws.Range("A6").Value = "ulkemis66"
$ws.Range("A7").Value = "ulkemis77"
$ws.Range("A8").Value = "ulkemis88"

# Column B (rows 1-8): "umisN" -> "uisN1"
$ws.Range("B1").Value = "uis11"
$ws.Range("B2").Value = "uis21"
$ws.Range("B3").Value = "uis31"
$ws.Range("B4").Value = "uis41"
$ws.Range("B5").Value = "uis51"
$ws.Range("B6").Value = "uis61"
$ws.Range("B7").Value = "uis71"
$ws.Range("B8").Value = "uis81"

# Update the selection on this sheet to match the new cursor position
$ws.Range("B9:B10").Select() | Out-Null
